$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.595.77'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +3.59%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.603.22'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.72%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.00'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.515'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '26.85'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +8.13%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.27'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -2.10%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.251'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0599'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.14%  '
$ws.Range("E12").Value = '  +1.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.832.65'
$ws.Range("D13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.602.01'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +2.63%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.612.28'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +3.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.537'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +4.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.72'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.61'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +3.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.40'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +5.23%  '
$ws.Range("E20").Value = '  +4.30%  '
$ws.Range("E21").Value = '  +1.07%  '
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("E23").Value = '  +1.56%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.24'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.89%  '
$ws.Range("E25").Value = '  +0.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.12'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +1.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.35'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +3.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.109'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.60%  '
$ws.Range("E29").Value = '  +2.51%  '
$ws.Range("E30").Value = '  +0.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0476'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +3.85%  '
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("E33").Value = '  +1.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.438.25'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.14'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +4.58%  '
$ws.Range("E36").Value = '  +4.89%  '
$ws.Range("E37").Value = '  -0.94%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.81'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +3.08%  '
$ws.Range("E39").Value = '  +0.73%  '
$ws.Range("E40").Value = '  +2.62%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.538'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +4.27%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.95'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +0.48%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0491'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +6.28%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '53.22'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +28.43%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.799'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.78%  '
$ws.Range("E46").Value = '  +0.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +21.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '65.81'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +3.92%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.30'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.742.43'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.69%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '86.41'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.46%  '
